$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1717.0952
$ws.Range("I112").Value = 766.6667
$ws.Range("J112").Value = 1875.5
$ws.Range("K112").Value = 2300.0001
$ws.Range("L112").Value = 5626.5
$ws.Range("M112").Value = -1192.0001
$ws.Range("N112").Value = -7842.5
$ws.Range("H113").Value = 2669.6155
$ws.Range("I113").Value = 1681
$ws.Range("J113").Value = 3287.5
$ws.Range("K113").Value = 1681
$ws.Range("L113").Value = 3287.5
$ws.Range("M113").Value = 1573
$ws.Range("N113").Value = -9795.5
$ws.Range("H116").Value = 2659.1177
$ws.Range("I116").Value = 1400
$ws.Range("J116").Value = 2827
$ws.Range("K116").Value = 1400
$ws.Range("L116").Value = 2827
$ws.Range("M116").Value = 2042
$ws.Range("N116").Value = -9711
$ws.Range("H132").Value = 4764967
$ws.Range("I132").Value = 5497692.5
$ws.Range("J132").Value = 2251.5
$ws.Range("K132").Value = 16493077.5
$ws.Range("L132").Value = 6754.5
$ws.Range("M132").Value = -16490547.5
$ws.Range("N132").Value = -11814.5
$ws.Range("H137").Value = 26317508
$ws.Range("I137").Value = 1360.3158
$ws.Range("K137").Value = 4080.9474
$ws.Range("M137").Value = -1530.9474
$ws.Range("H138").Value = 2612.7937
$ws.Range("I138").Value = 1263.9642
$ws.Range("K138").Value = 3791.8926
$ws.Range("M138").Value = 1348.1074

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1069.7391
$ws.Range("I2").Value = 608.2727
$ws.Range("J2").Value = 1492.75
$ws.Range("K2").Value = 608.2727
$ws.Range("L2").Value = 1492.75
$ws.Range("M2").Value = -495.2727
$ws.Range("N2").Value = -1718.75
$ws.Range("H61").Value = 1911.8572
$ws.Range("I61").Value = 1647.1666
$ws.Range("K61").Value = 1647.1666
$ws.Range("M61").Value = -1435.1666
$ws.Range("H63").Value = 1115370.9
$ws.Range("I63").Value = 1669056.4
$ws.Range("K63").Value = 1669056.4
$ws.Range("M63").Value = -1668370.4
$ws.Range("H66").Value = 1115370.9
$ws.Range("I66").Value = 1669056.4
$ws.Range("K66").Value = 8345282
$ws.Range("M66").Value = -8341850
$ws.Range("H116").Value = 1069.7391
$ws.Range("I116").Value = 608.2727
$ws.Range("J116").Value = 1492.75
$ws.Range("K116").Value = 608.2727
$ws.Range("L116").Value = 1492.75
$ws.Range("M116").Value = 1685.7273
$ws.Range("N116").Value = -6080.75
$ws.Range("H136").Value = 1911.8572
$ws.Range("I136").Value = 1647.1666
$ws.Range("K136").Value = 4941.4998
$ws.Range("M136").Value = -2391.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1069.7391
$ws.Range("I3").Value = 608.2727
$ws.Range("J3").Value = 1492.75
$ws.Range("K3").Value = 608.2727
$ws.Range("L3").Value = 1492.75
$ws.Range("M3").Value = -494.2727
$ws.Range("N3").Value = -1720.75
$ws.Range("H113").Value = 100040
$ws.Range("I113").Value = 100040
$ws.Range("K113").Value = 100040
$ws.Range("M113").Value = -97870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14710305
$ws.Range("I31").Value = 2063.8857
$ws.Range("K31").Value = 2063.8857
$ws.Range("M31").Value = -1768.8857
$ws.Range("H34").Value = 14710305
$ws.Range("I34").Value = 2063.8857
$ws.Range("K34").Value = 2063.8857
$ws.Range("M34").Value = -1861.8857
$ws.Range("H58").Value = 7408983.5
$ws.Range("I58").Value = 1431.3715
$ws.Range("J58").Value = 33335416
$ws.Range("K58").Value = 1431.3715
$ws.Range("L58").Value = 33335416
$ws.Range("M58").Value = -1228.3715
$ws.Range("N58").Value = -33335822
$ws.Range("H99").Value = 2344
$ws.Range("I99").Value = 2054.1667
$ws.Range("J99").Value = 2611.5386
$ws.Range("K99").Value = 2054.1667
$ws.Range("L99").Value = 2611.5386
$ws.Range("M99").Value = -556.1667000000002
$ws.Range("N99").Value = -5607.5386
$ws.Range("H107").Value = 902.9394
$ws.Range("I107").Value = 912.8182
$ws.Range("J107").Value = 883.1818
$ws.Range("K107").Value = 912.8182
$ws.Range("L107").Value = 883.1818
$ws.Range("M107").Value = 1007.1818
$ws.Range("N107").Value = -4723.1818
$ws.Range("H126").Value = 2344
$ws.Range("I126").Value = 2054.1667
$ws.Range("J126").Value = 2611.5386
$ws.Range("K126").Value = 6162.500100000001
$ws.Range("L126").Value = 7834.6158
$ws.Range("M126").Value = -3692.500100000001
$ws.Range("N126").Value = -12774.6158
$ws.Range("H134").Value = 4071.5312
$ws.Range("I134").Value = 4192.967
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 12578.901
$ws.Range("L134").Value = 6750
$ws.Range("M134").Value = -10043.901
$ws.Range("N134").Value = -11820
$ws.Range("H136").Value = 7408983.5
$ws.Range("I136").Value = 1431.3715
$ws.Range("J136").Value = 33335416
$ws.Range("K136").Value = 4294.1145
$ws.Range("L136").Value = 100006248
$ws.Range("M136").Value = -1744.1145
$ws.Range("N136").Value = -100011348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 707.65955
$ws.Range("I5").Value = 309.0645
$ws.Range("J5").Value = 1479.9375
$ws.Range("K5").Value = 927.1935000000001
$ws.Range("L5").Value = 4439.8125
$ws.Range("M5").Value = -815.1935000000001
$ws.Range("N5").Value = -4663.8125
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
$ws.Range("H68").Value = 552.8570999999999
$ws.Range("I68").Value = 545
$ws.Range("J68").Value = 600
$ws.Range("K68").Value = 1635
$ws.Range("L68").Value = 1800
$ws.Range("M68").Value = -824
$ws.Range("N68").Value = -3422
$ws.Range("H71").Value = 552.8570999999999
$ws.Range("I71").Value = 545
$ws.Range("J71").Value = 600
$ws.Range("K71").Value = 4905
$ws.Range("L71").Value = 5400
$ws.Range("M71").Value = -849
$ws.Range("N71").Value = -13512
$ws.Range("H131").Value = 776.09
$ws.Range("I131").Value = 421.66666
$ws.Range("J131").Value = 798.71277
$ws.Range("K131").Value = 1264.99998
$ws.Range("L131").Value = 2396.13831
$ws.Range("M131").Value = 3775.00002
$ws.Range("N131").Value = -12476.13831
$ws.Range("H135").Value = 707.65955
$ws.Range("I135").Value = 309.0645
$ws.Range("J135").Value = 1479.9375
$ws.Range("K135").Value = 2781.5805
$ws.Range("L135").Value = 13319.4375
$ws.Range("M135").Value = -246.5805
$ws.Range("N135").Value = -18389.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 7000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 7000
$ws.Range("K47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -8136
$ws.Range("H69").Value = 44000
$ws.Range("J69").Value = 44000
$ws.Range("L69").Value = 44000
$ws.Range("N69").Value = -45498
$ws.Range("H72").Value = 44000
$ws.Range("J72").Value = 44000
$ws.Range("L72").Value = 132000
$ws.Range("N72").Value = -139488
$ws.Range("H123").Value = 30200
$ws.Range("J123").Value = 30200
$ws.Range("L123").Value = 30200
$ws.Range("N123").Value = -35100
$ws.Range("H132").Value = 10003.533
$ws.Range("I132").Value = 11859.728
$ws.Range("J132").Value = 4899
$ws.Range("K132").Value = 35579.18399999999
$ws.Range("L132").Value = 14697
$ws.Range("M132").Value = -33049.18399999999
$ws.Range("N132").Value = -19757

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 29413652
$ws.Range("I7").Value = 1552.1818
$ws.Range("K7").Value = 1552.1818
$ws.Range("M7").Value = -1440.1818
$ws.Range("H40").Value = 2669.261
$ws.Range("I40").Value = 2394.65
$ws.Range("K40").Value = 2394.65
$ws.Range("M40").Value = -2258.65
$ws.Range("H122").Value = 10752
$ws.Range("I122").Value = 17004
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 51012
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -48562
$ws.Range("N122").Value = -18400
$ws.Range("H126").Value = 29413652
$ws.Range("I126").Value = 1552.1818
$ws.Range("K126").Value = 4656.5454
$ws.Range("M126").Value = -2186.5454
$ws.Range("H127").Value = 61900
$ws.Range("J127").Value = 61900
$ws.Range("L127").Value = 61900
$ws.Range("N127").Value = -71820
$ws.Range("H132").Value = 7206.25
$ws.Range("I132").Value = 8589.888999999999
$ws.Range("J132").Value = 3055.3333
$ws.Range("K132").Value = 25769.667
$ws.Range("L132").Value = 9165.999899999999
$ws.Range("M132").Value = -23239.667
$ws.Range("N132").Value = -14225.9999
$ws.Range("H133").Value = 22608.666
$ws.Range("J133").Value = 22608.666
$ws.Range("L133").Value = 22608.666
$ws.Range("N133").Value = -27668.666
